# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gains a new (blank) column before the
# existing "Late" column, to make room for a variable-instalment related
# column in this loan repayment schedule template. The sheet also becomes
# the active/selected sheet & cell in the workbook (as the author left it
# after editing the "Repayment schedule" tab), replacing "NewLoanInput"
# which was previously selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (pushing old N/O/P -> O/P/Q), mirroring the
# new empty "N" header/data column that appears ahead of "Late"/"Outstanding".
$ws.Columns("N").Insert()

# Match the new column's width to its neighbour (column M) rather than
# leaving it at the sheet default.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection on this sheet at R9 (just past the last used column),
# which also makes "Repayment schedule" the active tab of the workbook.
$ws.Range("R9").Select()
